$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 44699
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 56
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("S5").Value = 1200
$ws.Range("D6").Value = 44699
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("S6").Value = 1000
$ws.Range("D7").Value = 44306
$ws.Range("M7").Value = 45
$ws.Range("D8").Value = 44315
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("S8").Value = 1000
$ws.Range("D9").Value = 44312
$ws.Range("M9").Value = 48
$ws.Range("D10").Value = 44321
$ws.Range("M10").Value = 58
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 9000
$ws.Range("S10").Value = 900
$ws.Range("D11").Value = 44314
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 47
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 9000
$ws.Range("S11").Value = 900
$ws.Range("D12").Value = 44326
$ws.Range("M12").Value = 65
$ws.Range("D13").Value = 44326
$ws.Range("M13").Value = 67
$ws.Range("N13").Value = 8000
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 8000
$ws.Range("S13").Value = 800
$ws.Range("D14").Value = 44323
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 10000
$ws.Range("S14").Value = 1000
$ws.Range("D15").Value = 44323
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 9000
$ws.Range("O15").Value = 9000
$ws.Range("P15").Value = 9000
$ws.Range("S15").Value = 900
$ws.Range("D16").Value = 44329
$ws.Range("L16").Value = 'Primera'
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 9000
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 900
$ws.Range("D17").Value = 44329
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 800
$ws.Range("D18").Value = 44302
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 45
$ws.Range("R18").Value = 'Provincia de Quillota'
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 47
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 10000
$ws.Range("P19").Value = 10000
$ws.Range("S19").Value = 1000
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = 9000
$ws.Range("O20").Value = 9000
$ws.Range("P20").Value = 9000
$ws.Range("S20").Value = 900
$ws.Range("D21").Value = 44343
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 58
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("R21").Value = 'Región Metropolitana'
$ws.Range("S21").Value = 800
$ws.Range("D23").Value = 44308
$ws.Range("M23").Value = 45
$ws.Range("D24").Value = 44308
$ws.Range("M24").Value = 48
$ws.Range("D25").Value = 44301
$ws.Range("M25").Value = 45
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 10000
$ws.Range("S25").Value = 1000
$ws.Range("D26").Value = 44328
$ws.Range("N26").Value = 8000
$ws.Range("O26").Value = 8000
$ws.Range("P26").Value = 8000
$ws.Range("S26").Value = 800
$ws.Range("D27").Value = 44328
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 48
$ws.Range("N27").Value = 7000
$ws.Range("O27").Value = 7000
$ws.Range("P27").Value = 7000
$ws.Range("S27").Value = 700
$ws.Range("D28").Value = 44319
$ws.Range("M28").Value = 68
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("R28").Value = 'Provincia de Quillota'
$ws.Range("S28").Value = 1000
$ws.Range("D29").Value = 44319
$ws.Range("M29").Value = 57
$ws.Range("R29").Value = 'Provincia de Quillota'
$ws.Range("D30").Value = 44322
$ws.Range("M30").Value = 56
$ws.Range("D31").Value = 44322
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 40
$ws.Range("N31").Value = 8000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 8000
$ws.Range("S31").Value = 800
$ws.Range("D32").Value = 44307
$ws.Range("M32").Value = 40
